$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.284.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.895.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("E11").Value = "  -5.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.518.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.908.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("E17").Value = "  +6.79%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.278.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "425.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -5.20%  "
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "687.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.127"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "68.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.88%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.430"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0483"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.27%  "
$ws.Range("E48").Value = "  +6.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.757.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0344"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
